# Regenerate the "K" column (column G) values for the save_data sheet.
# Per commit message: "regen save_data to use K instead of Strike#,
# regen std/mean, calc and write s_vals" -- the recalculated K values
# are written back into column G for each data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 2
    3  = 2
    4  = 5
    5  = 1
    6  = 1
    7  = 1
    8  = 1
    9  = 0
    10 = 0
    11 = 2
    12 = 6
    13 = 2
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
